$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Update the date in A1 (month moved forward: 2024-04-24 -> 2024-05-24)
$ws.Range("A1").Value = 45436

# Update prices
$ws.Range("D28").Value = 47157
$ws.Range("D29").Value = 31460

# Re-create the merged cell regions so they are re-serialized in the
# expected order (unmerge then merge again in the desired sequence).
$ws.Range("A1:D1").UnMerge()
$ws.Range("A9:D9").UnMerge()
$ws.Range("B29:C29").UnMerge()
$ws.Range("B28:C28").UnMerge()
$ws.Range("B27:C27").UnMerge()
$ws.Range("A10:D10").UnMerge()

$ws.Range("A1:D1").Merge()
$ws.Range("A9:D9").Merge()
$ws.Range("B29:C29").Merge()
$ws.Range("B28:C28").Merge()
$ws.Range("B27:C27").Merge()
$ws.Range("A10:D10").Merge()
